$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{}
$data[2] = @(1.02,1.02020935029036,1.028186279667556,1.021220136258266,1.035130770811575,1.023594999628091,1.025407753437642,1.031003829276111,1.024058110239224,1.037928256868433,1.012452158818634)
$data[3] = @(1.02,1.02132920632262,1.029220324904256,1.02217679549443,1.036329240208038,1.023504579208683,1.026163808921131,1.031845146685901,1.02482076425516,1.038935022641682,1.012711453180477)
$data[4] = @(1.02,1.022054374107675,1.02989025789283,1.022796603130266,1.037105892393236,1.023443335729026,1.026653067329835,1.032389772332395,1.025314422923699,1.039587063615902,1.012879013867214)
$data[5] = @(1.02,1.022359366053031,1.03017209807532,1.023057358075439,1.037432675626427,1.023416932628374,1.026858761222235,1.032618789809561,1.02552199781915,1.03986132463169,1.012949403544407)
$data[6] = @(1.02,1.022410583268608,1.030219432033327,1.023101150987447,1.03748756032105,1.023412460900047,1.026893298700969,1.032657246189084,1.025556852948643,1.039907382655753,1.012961219178527)
$data[7] = @(1.02,1.022058448907757,1.029894023067784,1.022800086612803,1.037110257791373,1.023442985510381,1.026655815783625,1.032392832251586,1.025317196389367,1.03959072774283,1.012879954625311)
$data[8] = @(1.02,1.020587698224991,1.028535566922472,1.02154328106708,1.035535558714178,1.023565007130014,1.02566325770935,1.031288107421598,1.024315817406213,1.038268374782371,1.012539834050499)
$data[9] = @(1.02,1.018000199174954,1.026148197175385,1.019334648703178,1.032769622552555,1.023759139049229,1.023914545290446,1.029343251318956,1.022552566921314,1.035942774372515,1.01193881860297)
$data[10] = @(1.02,1.01627794640003,1.024560915205688,1.017866280673537,1.030931612079897,1.023874603872166,1.022748931942025,1.02804788840283,1.021377945132637,1.034395417686236,1.011537020073727)
$data[11] = @(1.02,1.015532832793117,1.023874621531023,1.017231423147907,1.030137137101096,1.023921309042974,1.022244252338499,1.027487265712157,1.020869528558607,1.033726112216236,1.011362772239941)
$data[12] = @(1.02,1.015256158558723,1.02361985292866,1.016995752314614,1.029842242742463,1.023938164268129,1.022056797388697,1.027279066933559,1.020680710271985,1.033477608559682,1.011298008889888)
$data[13] = @(1.02,1.015315501877299,1.023674494811126,1.017046298014399,1.029905489135307,1.023934571058996,1.022097006866215,1.027323724407619,1.020721211068086,1.033530908639985,1.011311902654673)
$data[14] = @(1.02,1.015509960899997,1.023853559191552,1.017211939595774,1.030112756781235,1.023922712353276,1.022228757142771,1.027470055092403,1.020853920163529,1.033705568662979,1.011357419695646)
$data[15] = @(1.02,1.015629785966005,1.023963906722703,1.017314015857491,1.030240488896373,1.023915340507297,1.022309933560281,1.027560219740042,1.020935690623516,1.033813196524735,1.011385458955995)
$data[16] = @(1.02,1.016327410391852,1.02460648350824,1.017908434205112,1.030984368116399,1.023871434968932,1.022782426655572,1.028085100914511,1.021411691318446,1.034439852205829,1.011548578722276)
$data[17] = @(1.02,1.016765180585314,1.025009825113517,1.018281552923018,1.031451357092457,1.023843013862435,1.02307881924232,1.028414419212291,1.021710328067462,1.034833126539841,1.011650828170646)
$data[18] = @(1.02,1.017020585694692,1.025245184907689,1.018499278954748,1.031723878568846,1.023826118601784,1.023251703773304,1.028606531847662,1.021884537369873,1.035062585361173,1.01171044283295)
$data[19] = @(1.02,1.01710768264132,1.02532545305531,1.018573533544294,1.031816824218712,1.023820303822167,1.023310653575343,1.028672041854562,1.02194394155865,1.035140836533132,1.011730765528751)
$data[20] = @(1.02,1.016718205667347,1.024966540297781,1.018241511279632,1.031401239666091,1.023846096030192,1.023047018736898,1.028379083701167,1.021678285150478,1.034790924838149,1.011639860431314)
$data[21] = @(1.02,1.015452694973948,1.02380082501687,1.017163158317031,1.03005171586077,1.02392621804881,1.022189959819831,1.027426963211209,1.020814839808249,1.033654132711969,1.011344017163625)
$data[22] = @(1.02,1.014657562323476,1.023068769509992,1.016485985348941,1.029204425698047,1.023973741382756,1.021651124695448,1.026828566936165,1.020272132494018,1.032940000555417,1.011157777683047)
$data[23] = @(1.02,1.015079025937756,1.023456763053269,1.016844888856979,1.029653475863672,1.023948818288725,1.021936768466199,1.0271457655103,1.020559815364791,1.03331851747302,1.011256528631875)
$data[24] = @(1.02,1.016739431421464,1.024986098541521,1.018259604089706,1.031423885157301,1.023844704313072,1.023061388004713,1.028395050214526,1.021692763904882,1.034809993758543,1.01164481636003)
$data[25] = @(1.02,1.018668642350705,1.026764631418412,1.019904919118999,1.03348363322587,1.02371141742819,1.02436659444327,1.029845829691816,1.023008254574207,1.036543460030294,1.012094393749036)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Host "Updated vm_pu values for case with 380 kV"
